$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

function Add-LogRow {
    param(
        [object]$ws,
        [int]$RowNum,
        [int]$Template,
        [object[]]$Values
    )

    # Copy the whole template row (formatting + values) into the destination
    # row first so the new row inherits borders/fill/number-format, then
    # overwrite the six cell values explicitly.
    $srcRange = $ws.Range("A$Template`:F$Template")
    $dstRange = $ws.Range("A$RowNum`:F$RowNum")
    $srcRange.Copy($dstRange)

    # The template's row height isn't part of the A:F cell copy, so carry it
    # over explicitly (templates are 15/30/45 depending on wrapped text).
    # Leave rows at the sheet's default (15) alone so they don't pick up a
    # spurious explicit height.
    $tmplHeight = $ws.Rows.Item($Template).RowHeight
    if ($tmplHeight -ne 15) {
        $ws.Rows.Item($RowNum).RowHeight = $tmplHeight
    }

    $cols = @("A","B","C","D","E","F")
    for ($i = 0; $i -lt 6; $i++) {
        $col = $cols[$i]
        $val = $Values[$i]
        $cell = $ws.Range("$col$RowNum")
        if ([string]::IsNullOrEmpty([string]$val)) {
            $cell.ClearContents()
        } elseif ($col -eq "B" -and $Template -ne 5) {
            $cell.Value = [double]$val
        } else {
            $cell.Value = $val
        }
    }
}

# --- Fix the lockup time on the existing Wed rows: 2200 -> 2150 ("Bug when
# removing MC"). Style/formatting for these cells is untouched.
foreach ($r in 289..292) {
    $ws.Range("C$r").Value = "2150"
}

# --- Append the Thu/Fri/Sat/Mon log entries.
Add-LogRow $ws 296 5 @('', 'THURSDAY', '', '', '', '')
Add-LogRow $ws 297 20 @('AV Shutdown', '42642', '1800', 'BC', '320', 'LEAVE ALL EQUIPMENT IN ROOM. JUST LOG OFF PC AND PROJECTOR AND LOCK IN. Key for room in CB 121A storeroom.')
Add-LogRow $ws 298 2 @('Other', '42642', '1800', 'BC', '320', 'Remote for projector is on PC cart - please leave there.')
Add-LogRow $ws 299 2 @('AV Shutdown', '42642', '1730', 'LSB', '101', 'Make sure neck mic goes back to drawer and log off touchscreen.')
Add-LogRow $ws 300 2 @('AV Shutdown', '42642', '1730', 'LSB', '107', 'Make sure neck mic goes back to drawer and log off touchscreen.')
Add-LogRow $ws 301 2 @('AV Shutdown', '42642', '1900', 'LSB', '103', 'Make sure neck mic goes back to drawer and log off touchscreen.')
Add-LogRow $ws 302 2 @('AV Shutdown', '42642', '1900', 'LSB', '105', 'Make sure neck mic goes back to drawer and log off touchscreen.')
Add-LogRow $ws 303 20 @('Demo', '42642', '1800', 'CLH', 'J', 'NO CEILING PROJECTOR - Use roll in PC and Projector that is in room. Make sure client is okay.')
Add-LogRow $ws 304 20 @('Demo', '42642', '1800', 'CLH', 'M', 'NO CEILING PROJECTOR - Use roll in PC and Projector that is in room. Make sure client is okay.')
Add-LogRow $ws 305 20 @('Demo', '42642', '1900', 'CLH', 'J', 'NO CEILING PROJECTOR - Use roll in PC and Projector that is in room. Make sure client is okay.')
Add-LogRow $ws 306 20 @('Demo', '42642', '1900', 'CLH', 'M', 'NO CEILING PROJECTOR - Use roll in PC and Projector that is in room. Make sure client is okay.')
Add-LogRow $ws 307 20 @('Lockup', '42642', '1730', 'CLH', 'K', 'PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS.')
Add-LogRow $ws 308 20 @('Lockup', '42642', '2150', 'CLH', 'H', 'PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS.')
Add-LogRow $ws 309 6 @('Lockup', '42642', '2150', 'CLH', 'M', 'LEAVE ROLL IN PC AND PROJECTOR IN ROOM - JUST TURN OFF. PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lasonde 1011 office. PLEASE LOCK ALL 4 DOORS.')
Add-LogRow $ws 310 6 @('Lockup', '42642', '2150', 'CLH', 'J', 'LEAVE ROLL IN PC AND PROJECTOR IN ROOM - JUST TURN OFF. PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lasonde 1011 office. PLEASE LOCK ALL 4 DOORS.')
Add-LogRow $ws 315 5 @('', 'FRIDAY', '', '', '', '')
Add-LogRow $ws 316 2 @('AV Shutdown', '42643', '1530', 'LSB', '101', 'Make sure neck mic goes back to drawer and log off touchscreen.')
Add-LogRow $ws 317 2 @('AV Shutdown', '42643', '1530', 'LSB', '103', 'Make sure neck mic goes back to drawer and log off touchscreen.')
Add-LogRow $ws 318 2 @('AV Shutdown', '42643', '1530', 'LSB', '105', 'Make sure neck mic goes back to drawer and log off touchscreen.')
Add-LogRow $ws 319 2 @('AV Shutdown', '42643', '1530', 'LSB', '106', 'Make sure neck mic goes back to drawer and log off touchscreen.')
Add-LogRow $ws 320 2 @('AV Shutdown', '42643', '1530', 'LSB', '107', 'Make sure neck mic goes back to drawer and log off touchscreen.')
Add-LogRow $ws 321 20 @('AV Shutdown', '42643', '1530', 'BC', '320', 'LEAVE ALL EQUIPMENT IN ROOM. JUST LOG OFF PC AND PROJECTOR AND LOCK IN. Key for room in CB 121A storeroom.')
Add-LogRow $ws 322 20 @('Lockup', '42643', '1730', 'CLH', 'H', 'PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS.')
Add-LogRow $ws 326 5 @('', 'SATURDAY', '', '', '', '')
Add-LogRow $ws 327 20 @('Lockup', '42644', '1300', 'CLH', 'K', 'PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS.')
Add-LogRow $ws 328 6 @('Lockup', '42644', '1300', 'CLH', 'J', 'LEAVE ROLL IN PC AND PROJECTOR IN ROOM - JUST TURN OFF. PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lasonde 1011 office. PLEASE LOCK ALL 4 DOORS.')
Add-LogRow $ws 329 6 @('Lockup', '42644', '1330', 'CLH', 'M', 'LEAVE ROLL IN PC AND PROJECTOR IN ROOM - JUST TURN OFF. PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lasonde 1011 office. PLEASE LOCK ALL 4 DOORS.')
Add-LogRow $ws 330 20 @('Lockup', '42644', '1400', 'CLH', 'H', 'PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS.')
Add-LogRow $ws 334 5 @('', 'MONDAY', '', '', '', '')
Add-LogRow $ws 335 2 @('AV Shutdown', '42646', '1630', 'LSB', '101', 'Make sure neck mic goes back to drawer and log off touchscreen.')
Add-LogRow $ws 336 2 @('AV Shutdown', '42646', '1900', 'LSB', '103', 'Make sure neck mic goes back to drawer and log off touchscreen.')
Add-LogRow $ws 337 2 @('AV Shutdown', '42646', '1900', 'LSB', '105', 'Make sure neck mic goes back to drawer and log off touchscreen.')
Add-LogRow $ws 338 2 @('AV Shutdown', '42646', '1730', 'LSB', '106', 'Make sure neck mic goes back to drawer and log off touchscreen.')
Add-LogRow $ws 339 2 @('AV Shutdown', '42646', '1900', 'LSB', '107', 'Make sure neck mic goes back to drawer and log off touchscreen.')
Add-LogRow $ws 340 20 @('Demo', '42646', '1900', 'CLH', 'J', 'NO CEILING PROJECTOR - Use roll in PC and Projector that is in room. Make sure client is okay.')
Add-LogRow $ws 341 6 @('Lockup', '42646', '1900', 'CLH', 'M', 'LEAVE ROLL IN PC AND PROJECTOR IN ROOM - JUST TURN OFF. PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lasonde 1011 office. PLEASE LOCK ALL 4 DOORS.')
Add-LogRow $ws 342 20 @('Lockup', '42646', '2100', 'CLH', 'H', 'PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS.')
Add-LogRow $ws 343 20 @('Lockup', '42646', '2150', 'CLH', 'K', 'PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lassonde 1011 office. PLEASE LOCK ALL 4 DOORS.')
Add-LogRow $ws 344 6 @('Lockup', '42646', '2150', 'CLH', 'J', 'LEAVE ROLL IN PC AND PROJECTOR IN ROOM - JUST TURN OFF. PLEASE LOCK ROOM. ALLEN KEY (with yellow handle) and CLH key is on keyrack in Lasonde 1011 office. PLEASE LOCK ALL 4 DOORS.')

# --- Match the recorded scroll position / selection from the diff.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 329
$ws.Range("A344").Select()
